$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for B1 (unicode micro sign mu, U+03BC)
$ws.Range("B1").Value = "Niacine par semaine [μg]"

# Update column B values for rows 2-79
$ws.Range("B2").Value = 177981.225
$ws.Range("B3").Value = 94381.60000000001
$ws.Range("B4").Value = 67700.52499999999
$ws.Range("B5").Value = 71081.57500000001
$ws.Range("B6").Value = 137804.5
$ws.Range("B7").Value = 82495.575
$ws.Range("B8").Value = 126552.5
$ws.Range("B9").Value = 90354.59999999999
$ws.Range("B10").Value = 205466.375
$ws.Range("B11").Value = 137023.25
$ws.Range("B12").Value = 93649.575
$ws.Range("B13").Value = 106441.325
$ws.Range("B14").Value = 89049.35000000001
$ws.Range("B15").Value = 36431.10000000001
$ws.Range("B16").Value = 141812.125
$ws.Range("B17").Value = 140692.5
$ws.Range("B18").Value = 140565.15
$ws.Range("B19").Value = 205072.9
$ws.Range("B20").Value = 132757.75
$ws.Range("B21").Value = 50431.05
$ws.Range("B22").Value = 88933.47499999999
$ws.Range("B23").Value = 66639.95
$ws.Range("B24").Value = 147296.775
$ws.Range("B25").Value = 105460.825
$ws.Range("B26").Value = 62991.225
$ws.Range("B27").Value = 80606.5
$ws.Range("B28").Value = 75746.02499999999
$ws.Range("B29").Value = 122639.475
$ws.Range("B30").Value = 131828.9
$ws.Range("B31").Value = 130056.725
$ws.Range("B32").Value = 54337.125
$ws.Range("B33").Value = 100279.95
$ws.Range("B34").Value = 70846.47500000001
$ws.Range("B35").Value = 124616.225
$ws.Range("B36").Value = 39206.675
$ws.Range("B37").Value = 178707.775
$ws.Range("B38").Value = 188456.7
$ws.Range("B39").Value = 100736.925
$ws.Range("B40").Value = 66083.325
$ws.Range("B41").Value = 85475.97499999999
$ws.Range("B42").Value = 86002.47500000001
$ws.Range("B43").Value = 117164.875
$ws.Range("B44").Value = 103308.975
$ws.Range("B45").Value = 101879.75
$ws.Range("B46").Value = 99280.92499999999
$ws.Range("B47").Value = 120040.225
$ws.Range("B48").Value = 18657.625
$ws.Range("B49").Value = 150993.15
$ws.Range("B50").Value = 36602.4
$ws.Range("B51").Value = 73274.20000000001
$ws.Range("B52").Value = 92313.54999999999
$ws.Range("B53").Value = 48811.52499999999
$ws.Range("B54").Value = 120925.35
$ws.Range("B55").Value = 128656.25
$ws.Range("B56").Value = 64189.575
$ws.Range("B57").Value = 88626.35000000001
$ws.Range("B58").Value = 161753.65
$ws.Range("B59").Value = 62874.77499999999
$ws.Range("B60").Value = 39198.625
$ws.Range("B61").Value = 95701.52499999999
$ws.Range("B62").Value = 105590.65
$ws.Range("B63").Value = 69008.85000000001
$ws.Range("B64").Value = 42016.625
$ws.Range("B65").Value = 66011.02499999999
$ws.Range("B66").Value = 45632.925
$ws.Range("B67").Value = 99751.32500000001
$ws.Range("B68").Value = 61738.15
$ws.Range("B69").Value = 61422.6
$ws.Range("B70").Value = 49176.425
$ws.Range("B71").Value = 97283.92499999999
$ws.Range("B72").Value = 111093.85
$ws.Range("B73").Value = 62099.85
$ws.Range("B74").Value = 78525.375
$ws.Range("B75").Value = 111613.1
$ws.Range("B76").Value = 76328.72500000001
$ws.Range("B77").Value = 150102.075
$ws.Range("B78").Value = 81209
$ws.Range("B79").Value = 151451.4
